# Applies numeric corrections to the profit-calculation columns (H-N)
# across several job sheets, per the scheduled runner's refreshed pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 176.27272
$ws.Range("I5").Value = 149
$ws.Range("J5").Value = 299
$ws.Range("K5").Value = 149
$ws.Range("L5").Value = 299
$ws.Range("M5").Value = -34
$ws.Range("N5").Value = -529
$ws.Range("H9").Value = 188.27272
$ws.Range("I9").Value = 198
$ws.Range("J9").Value = 171.25
$ws.Range("K9").Value = 198
$ws.Range("L9").Value = 171.25
$ws.Range("M9").Value = -29
$ws.Range("N9").Value = -509.25
$ws.Range("H43").Value = 12086.692
$ws.Range("I43").Value = 29810.5
$ws.Range("J43").Value = 4209.4443
$ws.Range("K43").Value = 29810.5
$ws.Range("L43").Value = 4209.4443
$ws.Range("M43").Value = -29741.5
$ws.Range("N43").Value = -4347.4443
$ws.Range("H70").Value = 3184
$ws.Range("I70").Value = 3495
$ws.Range("J70").Value = 3059.6
$ws.Range("K70").Value = 10485
$ws.Range("L70").Value = 9178.799999999999
$ws.Range("M70").Value = -10215
$ws.Range("N70").Value = -9718.799999999999
$ws.Range("H73").Value = 3184
$ws.Range("I73").Value = 3495
$ws.Range("J73").Value = 3059.6
$ws.Range("K73").Value = 10485
$ws.Range("L73").Value = 9178.799999999999
$ws.Range("M73").Value = -9549
$ws.Range("N73").Value = -11050.8
$ws.Range("H141").Value = 5666.6665
$ws.Range("I141").Value = 5666.6665
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 16999.9995
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -11819.9995
$ws.Range("N141").Value = ""
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 270330.1
$ws.Range("I32").Value = 300662.72
$ws.Range("K32").Value = 300662.72
$ws.Range("M32").Value = -300375.72
$ws.Range("H61").Value = 5123.9165
$ws.Range("I61").Value = 1924.4667
$ws.Range("J61").Value = 10456.333
$ws.Range("K61").Value = 1924.4667
$ws.Range("L61").Value = 10456.333
$ws.Range("M61").Value = -1712.4667
$ws.Range("N61").Value = -10880.333
$ws.Range("H102").Value = 5765.6
$ws.Range("I102").Value = 4769.5
$ws.Range("K102").Value = 4769.5
$ws.Range("M102").Value = -3147.5
$ws.Range("H132").Value = 1790922
$ws.Range("I132").Value = 2277127.2
$ws.Range("J132").Value = 8169
$ws.Range("K132").Value = 6831381.600000001
$ws.Range("L132").Value = 24507
$ws.Range("M132").Value = -6828851.600000001
$ws.Range("N132").Value = -29567
$ws.Range("H136").Value = 5123.9165
$ws.Range("I136").Value = 1924.4667
$ws.Range("J136").Value = 10456.333
$ws.Range("K136").Value = 5773.4001
$ws.Range("L136").Value = 31368.999
$ws.Range("M136").Value = -3223.4001
$ws.Range("N136").Value = -36468.999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2316.2449
$ws.Range("I105").Value = 1688.5143
$ws.Range("J105").Value = 3885.5715
$ws.Range("K105").Value = 1688.5143
$ws.Range("L105").Value = 3885.5715
$ws.Range("M105").Value = 58.48569999999995
$ws.Range("N105").Value = -7379.5715
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1995.6
$ws.Range("I16").Value = 1995.6
$ws.Range("K16").Value = 1995.6
$ws.Range("M16").Value = -1708.6
$ws.Range("H48").Value = 10000
$ws.Range("I48").Value = 10000
$ws.Range("K48").Value = 10000
$ws.Range("M48").Value = -9524
$ws.Range("H99").Value = 8619.666999999999
$ws.Range("I99").Value = 8619.666999999999
$ws.Range("K99").Value = 8619.666999999999
$ws.Range("M99").Value = -7121.666999999999
$ws.Range("H105").Value = 16118.571
$ws.Range("I105").Value = 20686
$ws.Range("J105").Value = 4700
$ws.Range("K105").Value = 20686
$ws.Range("L105").Value = 4700
$ws.Range("M105").Value = -18939
$ws.Range("N105").Value = -8194
$ws.Range("H113").Value = 1995.6
$ws.Range("I113").Value = 1995.6
$ws.Range("K113").Value = 1995.6
$ws.Range("M113").Value = 174.4000000000001
$ws.Range("H126").Value = 8619.666999999999
$ws.Range("I126").Value = 8619.666999999999
$ws.Range("K126").Value = 25859.001
$ws.Range("M126").Value = -23389.001
$ws.Range("H132").Value = 2480.7334
$ws.Range("I132").Value = 2183.0833
$ws.Range("J132").Value = 3671.3333
$ws.Range("K132").Value = 6549.249899999999
$ws.Range("L132").Value = 11013.9999
$ws.Range("M132").Value = -4019.249899999999
$ws.Range("N132").Value = -16073.9999
$ws.Range("H141").Value = 349933
$ws.Range("J141").Value = 415319.6
$ws.Range("L141").Value = 415319.6
$ws.Range("N141").Value = -425679.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 326.25
$ws.Range("J23").Value = 133.25
$ws.Range("L23").Value = 399.75
$ws.Range("N23").Value = -869.75
$ws.Range("H62").Value = 1419.75
$ws.Range("I62").Value = 893
$ws.Range("K62").Value = 2679
$ws.Range("M62").Value = -1993
$ws.Range("H65").Value = 1419.75
$ws.Range("I65").Value = 893
$ws.Range("K65").Value = 8037
$ws.Range("M65").Value = -4605
$ws.Range("H132").Value = 861.63635
$ws.Range("I132").Value = 621.3333
$ws.Range("K132").Value = 5591.9997
$ws.Range("M132").Value = -3061.9997
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 61061
$ws.Range("I70").Value = 52886.8
$ws.Range("K70").Value = 52886.8
$ws.Range("M70").Value = -52616.8
$ws.Range("H73").Value = 61061
$ws.Range("I73").Value = 52886.8
$ws.Range("K73").Value = 52886.8
$ws.Range("M73").Value = -51950.8
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""
$ws.Range("H122").Value = 115312.78
$ws.Range("I122").Value = 170386.17
$ws.Range("J122").Value = 5166
$ws.Range("K122").Value = 511158.51
$ws.Range("L122").Value = 15498
$ws.Range("M122").Value = -508708.51
$ws.Range("N122").Value = -20398
$ws.Range("H132").Value = 11873.792
$ws.Range("I132").Value = 14530.143
$ws.Range("K132").Value = 43590.429
$ws.Range("M132").Value = -41060.429
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6276.125
$ws.Range("I7").Value = 4045
$ws.Range("J7").Value = 9994.666999999999
$ws.Range("K7").Value = 4045
$ws.Range("L7").Value = 9994.666999999999
$ws.Range("M7").Value = -3933
$ws.Range("N7").Value = -10218.667
$ws.Range("H126").Value = 6276.125
$ws.Range("I126").Value = 4045
$ws.Range("J126").Value = 9994.666999999999
$ws.Range("K126").Value = 12135
$ws.Range("L126").Value = 29984.001
$ws.Range("M126").Value = -9665
$ws.Range("N126").Value = -34924.001
$ws.Range("H132").Value = 2872
$ws.Range("I132").Value = 2756.3333
$ws.Range("J132").Value = 3045.5
$ws.Range("K132").Value = 8268.999899999999
$ws.Range("L132").Value = 9136.5
$ws.Range("M132").Value = -5738.999899999999
$ws.Range("N132").Value = -14196.5
$ws.Range("H136").Value = 9533.368
$ws.Range("I136").Value = 5133
$ws.Range("K136").Value = 15399
$ws.Range("M136").Value = -12849
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 25000
$ws.Range("J75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("N75").Value = -26872
$ws.Range("H78").Value = 25000
$ws.Range("J78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("N78").Value = -84360
$ws.Range("H132").Value = 3845.818
$ws.Range("I132").Value = 3438
$ws.Range("J132").Value = 4933.3335
$ws.Range("K132").Value = 10314
$ws.Range("L132").Value = 14800.0005
$ws.Range("M132").Value = -7784
$ws.Range("N132").Value = -19860.0005
